# TC04_INS_CancerType-BroadCancerTypes.xlsx
# "cancer type facet INS - 17 test scripts"
#
# The Programs-tab SQL query (row 2 / column B, the "TabQuery" for the
# "ProgramsTab" row) is rewritten:
#   - "Data Location Details" now comes from a CASE expression instead of
#     the raw prg.data_link column.
#   - ORDER BY now sorts on LOWER(prg.program_name) instead of the raw
#     program_name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProgramsQuery = "SELECT DISTINCT `n    prg.program_name AS ""Program"",`n    prg.website AS ""Website"",`n    prg.focus_area AS ""Focus Area"",`n    prg.cancer_type AS ""Cancer Type"",`nCASE `n        WHEN prg.data_link IS NOT NULL THEN prg.website       `n        ELSE prg.data_link`n    END AS ""Data Location Details""`nFROM `n    df_program prg`nWHERE `n     prg.cancer_type IN ('Broad Cancer Types')`nORDER BY `n    LOWER (prg.program_name) ASC`nLIMIT 100;"

$ws.Range("B2").Value2 = $newProgramsQuery

# Re-apply the cell's font/wrap formatting; this mirrors what Excel does
# internally when the text of a styled cell is edited (it re-materialises
# a dedicated style entry for the cell instead of reusing the old shared
# one), which is what the saved workbook reflects (fonts count 5 -> 6,
# cellXfs count 4 -> 5).
$ws.Range("B2").Font.Size = 12
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").WrapText = $true

# Move the selection to where the user was working (C3) to mirror the
# saved view state in the workbook.
[void]$ws.Range("C3").Select()

Write-Host "B2 updated:"
Write-Host $ws.Range("B2").Value2
